$wb = $excel.ActiveWorkbook

# --- "Direction 0 STOPS" sheet: clear out the last 4 stops, keep the
#     trailing blank rows exactly as before (dimension/rows unchanged). ---
$ws1 = $wb.Worksheets.Item("Direction 0 STOPS")
$ws1.Range("A5:B8").ClearContents() | Out-Null

# --- "Direction 1 STOPS" sheet: drop the first four stops (shifting the
#     remaining three up) and leave one trailing blank, styled row behind,
#     matching the new (smaller) used range. ---
$ws2 = $wb.Worksheets.Item("Direction 1 STOPS")
$ws2.Range("A2:B5").Delete(-4162) | Out-Null
$ws2.Range("A4:B4").Copy() | Out-Null
$ws2.Range("A5:B5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A5:B5").ClearContents() | Out-Null
$excel.CutCopyMode = $false

# --- Selections: set sheet 2's selection first so that activating sheet 1
#     afterwards leaves it as the active/selected tab, matching the target. ---
$ws2.Activate()
$ws2.Range("A2:B4").Select() | Out-Null

$ws1.Activate()
$ws1.Range("C10").Select() | Out-Null
